$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[0,80]"
$ws.Range("D7").Value = " 0-9; 10-19; 20-29; 30-39; 40-49; 50-59; 60-69; 70+"
$ws.Range("D15").Value = "[0,5]"
$ws.Range("D16").Value = "[1,13]"
$ws.Range("D45").Value = "[0,30]"
$ws.Range("D46").Value = "[0,30]"
$ws.Range("D51").Value = "[14,45]"
$ws.Range("D52").Value = "[0,12]"
$ws.Range("D55").Value = "[1,7]"
$ws.Range("D58").Value = "[0,6]"
$ws.Range("D59").Value = "[0,6]"
